$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (avoid Excel auto-numeric coercion
# stripping trailing zeros / reparsing dotted values as numbers or dates).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.230.71'
$ws.Range("E2").Value = '  +3.19%  '
$ws.Range("D3").Value = '1.582.86'
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = '212.60'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  +6.14%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").Value = '26.33'
$ws.Range("E8").Value = '  +10.44%  '
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").Value = '1.808.42'
$ws.Range("E12").Value = '  +1.89%  '
$ws.Range("D13").Value = '1.577.52'
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").Value = '29.252.60'
$ws.Range("E14").Value = '  +3.27%  '
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("E16").Value = '  +2.74%  '
$ws.Range("D17").Value = '62.85'
$ws.Range("E17").Value = '  +3.29%  '
$ws.Range("D18").Value = '238.64'
$ws.Range("E18").Value = '  +4.67%  '
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("E20").Value = '  +2.30%  '
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("D23").Value = '9.22'
$ws.Range("E23").Value = '  +3.27%  '
$ws.Range("E24").Value = '  +3.27%  '
$ws.Range("D25").Value = '154.32'
$ws.Range("E25").Value = '  +2.26%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '15.18'
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '0.108'
$ws.Range("E27").Value = '  +4.89%  '
$ws.Range("D28").Value = '6.37'
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").Value = '0.0471'
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("D33").Value = '1.427.59'
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("D34").Value = '3.09'
$ws.Range("E34").Value = '  +2.26%  '
$ws.Range("E35").Value = '  -3.12%  '
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("E37").Value = '  +9.05%  '
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("E40").Value = '  +3.70%  '
$ws.Range("E41").Value = '  +2.78%  '
$ws.Range("D42").Value = '54.00'
$ws.Range("E42").Value = '  +26.79%  '
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("E44").Value = '  +2.22%  '
$ws.Range("E45").Value = '  +2.36%  '
$ws.Range("D46").Value = '64.61'
$ws.Range("E46").Value = '  +4.32%  '
$ws.Range("D48").Value = '1.720.96'
$ws.Range("E48").Value = '  +2.07%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '85.70'
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '0.837'
$ws.Range("E50").Value = '  -5.93%  '
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -0.94%  '

Write-Output "Applied cryptos update."
